# Weekly update: push a new price observation into row 23 (Fruta / hortaliza,
# semanal) and shift the existing historical rows 23-61 down to 24-62,
# preserving the constant columns (A,B,C,E,F,G,I,N,O,Q,R) and moving only the
# per-record columns (D,H,J,K,L,M,P).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that vary per record and must shift down by one row.
$cols = @(4, 8, 10, 11, 12, 13, 16)   # D, H, J, K, L, M, P

# First, create the new last row (62) by copying what is currently in row 61
# (this becomes the bottom-most historical record after the shift).
foreach ($c in $cols) {
    $ws.Cells.Item(62, $c).Value = $ws.Cells.Item(61, $c).Value2
}
# Make sure the date cell in the new row carries the same date format as the
# rest of column D.
$ws.Cells.Item(62, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Shift rows 61 down to 24, working from the bottom up so we never overwrite
# a source row before it has been read.
for ($r = 60; $r -ge 23; $r--) {
    foreach ($c in $cols) {
        $ws.Cells.Item($r + 1, $c).Value = $ws.Cells.Item($r, $c).Value2
    }
}

# Fill in the constant columns for the newly-created row 62 (identical for
# every record in this sheet).
$ws.Cells.Item(62, 1).Value = 4
$ws.Cells.Item(62, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(62, 3).Value = "Los Lagos"
$ws.Cells.Item(62, 5).Value = 10
$ws.Cells.Item(62, 6).Value = 300000000
$ws.Cells.Item(62, 7).Value = "Espárragos"
$ws.Cells.Item(62, 9).Value = "Primera"
$ws.Cells.Item(62, 14).Value = "`$/kilo"
$ws.Cells.Item(62, 15).Value = "Provincia de Linares"
$ws.Cells.Item(62, 17).Value = 1
$ws.Cells.Item(62, 18).Value = "Hortaliza"

# Finally, write the new weekly record into row 23.
$ws.Cells.Item(23, 4).Value = 44880
$ws.Cells.Item(23, 8).Value = "Sin especificar"
$ws.Cells.Item(23, 10).Value = 600
$ws.Cells.Item(23, 11).Value = 1500
$ws.Cells.Item(23, 12).Value = 1700
$ws.Cells.Item(23, 13).Value = 1600
$ws.Cells.Item(23, 16).Value = 1600
